$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts for "Contenu du stage" (C#: 12 -> 11, ASSEMBLEUR: 1 -> 2)
$ws.Range("E16").Value = 11
$ws.Range("E19").Value = 2

# Update displayed percentages accordingly (force text, not numeric percent parsing)
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "73.33 %"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "13.33 %"

# Make sure the "Contenu du stage" pie chart (Chart 2) picks up the new
# figures from the worksheet.
$excel.CalculateFull()
$wb.RefreshAll()
$co = $ws.ChartObjects().Item(2)
$chart = $co.Chart
$chart.Refresh()
